$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.401.48"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.847.26"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6287"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07605"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2930"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.854.37"
$ws.Range("E12").Value = "  -6.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001088"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6793"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "2.098.68"
$ws.Range("E17").Value = "  -7.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.181"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "29.419.41"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.468"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.360"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.464"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.303"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05583"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.030"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "1.233.78"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.428"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9062"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.186"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4020"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.951"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.677"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -0.62%  "
